$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.883459091186523
$ws.Range("C2").Value = 5.528735637664795
$ws.Range("D2").Value = 13.01597785949707
$ws.Range("E2").Value = 57.85714340209961
